$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 690, shifting rows 690-737 down to 691-738.
$ws.Rows("690:690").Insert()

# Populate the newly inserted row 690 with the new price-report entry.
$ws.Range("A690").Value = 4
$ws.Range("B690").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C690").Value = "Los Lagos"
$ws.Range("D690").Value = 45021
$ws.Range("E690").Value = 10
$ws.Range("F690").Value = "Fruta"
$ws.Range("G690").Value = 100102
$ws.Range("H690").Value = "Cítricos"
$ws.Range("I690").Value = 100102005
$ws.Range("J690").Value = "Naranja"
$ws.Range("K690").Value = "Valencia"
$ws.Range("L690").Value = "Primera"
$ws.Range("M690").Value = 300
$ws.Range("N690").Value = 17000
$ws.Range("O690").Value = 18000
$ws.Range("P690").Value = 17500
$ws.Range("Q690").Value = '$/malla 18 kilos'
$ws.Range("R690").Value = "Región de O'Higgins"
$ws.Range("S690").Value = 972
$ws.Range("T690").Value = 18
